# "Generate Report for Archive" - refresh the localization status report:
# the zh-cn / de-de files have moved from "Ready for handoff" into
# "In Translation", and the (now shorter) Status column can be narrowed to fit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2)
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# Per-locale detail sheets: Status column (C2)
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# Narrow the Status columns to fit the shorter text
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
